$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from the (now-shifted) data columns into the new D:E columns
$ws.Range("F7:M102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = most recent quarter, E = prior quarter)
$ws.Range("D7").Value = 43464
$ws.Range("E7").Value = 43366
$ws.Range("D8").Value = 413000
$ws.Range("E8").Value = 408300
$ws.Range("D9").Value = 277800
$ws.Range("E9").Value = 280100
$ws.Range("D10").Value = 135200
$ws.Range("E10").Value = 128200
$ws.Range("D12").Value = 49200
$ws.Range("E12").Value = 46000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 6300
$ws.Range("E15").Value = 8500
$ws.Range("D17").Value = 405600
$ws.Range("E17").Value = 407700
$ws.Range("D18").Value = 7400
$ws.Range("E18").Value = 600
$ws.Range("D20").Value = -5400
$ws.Range("E20").Value = -9500
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = 30800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1900
$ws.Range("E23").Value = -9000
$ws.Range("D24").Value = 4400
$ws.Range("E24").Value = 2200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -2500
$ws.Range("E26").Value = -11100
$ws.Range("D27").Value = -2500
$ws.Range("E27").Value = -11100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 5400
$ws.Range("E32").Value = 9500
$ws.Range("D33").Value = -2500
$ws.Range("E33").Value = -11100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -2500
$ws.Range("E35").Value = -11100
$ws.Range("D38").Value = 43464
$ws.Range("E38").Value = 43366
$ws.Range("D41").Value = 356000
$ws.Range("E41").Value = 336300
$ws.Range("D42").Value = 367600
$ws.Range("E42").Value = 329200
$ws.Range("D43").Value = 193300
$ws.Range("E43").Value = 210400
$ws.Range("D44").Value = 313300
$ws.Range("E44").Value = 306400
$ws.Range("D45").Value = 45900
$ws.Range("E45").Value = 47500
$ws.Range("D46").Value = 1276200
$ws.Range("E46").Value = 1229800
$ws.Range("D47").Value = 48400
$ws.Range("E47").Value = 50200
$ws.Range("D48").Value = 675900
$ws.Range("E48").Value = 668300
$ws.Range("D49").Value = 994500
$ws.Range("E49").Value = 1001700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 19400
$ws.Range("E52").Value = 17600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3014600
$ws.Range("E54").Value = 2967600
$ws.Range("D57").Value = 143500
$ws.Range("E57").Value = 147200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 151300
$ws.Range("E59").Value = 146400
$ws.Range("D60").Value = 294700
$ws.Range("E60").Value = 293600
$ws.Range("D61").Value = 458000
$ws.Range("E61").Value = 452600
$ws.Range("D62").Value = 38300
$ws.Range("E62").Value = 22700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 795900
$ws.Range("E66").Value = 773800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -485900
$ws.Range("E72").Value = -483500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2218700
$ws.Range("E76").Value = 2193800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43464
$ws.Range("E80").Value = 43366
$ws.Range("D81").Value = -2500
$ws.Range("E81").Value = -11100
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = 39800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = "NA"
$ws.Range("E89").Value = 34000
$ws.Range("D91").Value = "NA"
$ws.Range("E91").Value = -36600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = "NA"
$ws.Range("E94").Value = -102300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = "NA"
$ws.Range("E100").Value = 285600
$ws.Range("D101").Value = "NA"
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = "NA"
$ws.Range("E102").Value = 217400

# Row 91 (Capital Expenditures) also had its historical F:J figures revised, not just shifted
$ws.Range("F91").Value = -57300
$ws.Range("G91").Value = -43200
$ws.Range("H91").Value = -48800
$ws.Range("I91").Value = -36500
$ws.Range("J91").Value = -30000
